# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.835.06"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "2.285.40"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "120.54"
$ws.Range("E5").Value = "  +7.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.51"
$ws.Range("E6").Value = "  +1.45%  "

$ws.Range("E7").Value = "  +5.29%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.636"
$ws.Range("E9").Value = "  +6.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.58"
$ws.Range("E10").Value = "  +2.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0954"
$ws.Range("E11").Value = "  +3.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.30"
$ws.Range("E12").Value = "  +7.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.62"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.922"
$ws.Range("E15").Value = "  +8.86%  "

$ws.Range("D16").Value = "2.628.96"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("D17").Value = "2.272.77"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "43.802.07"
$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("E19").Value = "  +4.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.98"
$ws.Range("E20").Value = "  +1.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.43"
$ws.Range("E21").Value = "  +2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").Value = "  +1.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.62"
$ws.Range("E23").Value = "  +3.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.67"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").Value = "  +2.60%  "

$ws.Range("E26").Value = "  +7.86%  "

$ws.Range("E27").Value = "  +1.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.79"
$ws.Range("E28").Value = "  +4.80%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.80"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.74"
$ws.Range("E32").Value = "  +2.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0932"
$ws.Range("E33").Value = "  +4.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  +5.68%  "

$ws.Range("E35").Value = "  +4.57%  "

$ws.Range("E36").Value = "  +14.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0388"
$ws.Range("E37").Value = "  +11.96%  "

$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("E39").Value = "  +4.77%  "

$ws.Range("E40").Value = "  +8.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.54"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.81"
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  +3.18%  "

$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("E45").Value = "  +2.80%  "

$ws.Range("E46").Value = "  -3.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "75.58"
$ws.Range("E47").Value = "  +47.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +4.01%  "

$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("E50").Value = "  +2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.53"
$ws.Range("E51").Value = "  +3.22%  "
